$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Gnas", "Gcgr", "ECs", 3, 1, 111.005114, 333.015342, 0.07617677061978854, 0.07617677061978853, 1, 0.3333333333333333, 0.2537453333333333, 0.761236, 0.1395119621451497, 0.1395119621451497, 28.16702965363467, 253.503266882712, 0.01062757073904769, 0.01062757073904769),
    @("ECs", "Gnas", "Gcgr", "FAPs", 3, 1, 111.005114, 333.015342, 0.07617677061978854, 0.07617677061978853, 3, 1, 1.194006666666667, 3.58202, 0.6564779367281227, 0.6564779367281226, 132.5408461500933, 1192.86761535084, 0.05000836920309026, 0.05000836920309024),
    @("ECs", "Gnas", "Gcgr", "M2", 3, 1, 111.005114, 333.015342, 0.07617677061978854, 0.07617677061978853, 1, 0.3333333333333333, 0.1154886666666667, 0.346466, 0.06349693324616998, 0.06349693324616996, 12.81983260904133, 115.378493481372, 0.004836991318953516, 0.004836991318953513),
    @("ECs", "Gnas", "Gcgr", "sCs", 3, 1, 111.005114, 333.015342, 0.07617677061978854, 0.07617677061978853, 3, 1, 0.2555663333333333, 0.766699, 0.1405131678805576, 0.1405131678805576, 28.36916996622867, 255.322529696058, 0.01070383935869708, 0.01070383935869707),
    @("FAPs", "Gnas", "Gcgr", "ECs", 3, 1, 930.1503093333332, 2790.450928, 0.6383115534299679, 0.6383115534299678, 1, 0.3333333333333333, 0.2537453333333333, 0.761236, 0.1395119621451497, 0.1395119621451497, 236.0213002918897, 2124.191702627008, 0.08905209727893341, 0.08905209727893337),
    @("FAPs", "Gnas", "Gcgr", "FAPs", 3, 1, 930.1503093333332, 2790.450928, 0.6383115534299679, 0.6383115534299678, 3, 1, 1.194006666666667, 3.58202, 0.6564779367281227, 0.6564779367281226, 1110.605670346062, 9995.451033114559, 0.4190374515854282, 0.419037451585428),
    @("FAPs", "Gnas", "Gcgr", "M2", 3, 1, 930.1503093333332, 2790.450928, 0.6383115534299679, 0.6383115534299678, 1, 0.3333333333333333, 0.1154886666666667, 0.346466, 0.06349693324616998, 0.06349693324616996, 107.4218190244942, 966.7963712204479, 0.04053082609840174, 0.04053082609840172),
    @("FAPs", "Gnas", "Gcgr", "sCs", 3, 1, 930.1503093333332, 2790.450928, 0.6383115534299679, 0.6383115534299678, 3, 1, 0.2555663333333333, 0.766699, 0.1405131678805576, 0.1405131678805576, 237.7151040051858, 2139.435936046672, 0.0896911784672046, 0.08969117846720458),
    @("M2", "Gnas", "Gcgr", "ECs", 3, 1, 72.47344466666667, 217.420334, 0.04973458223193756, 0.04973458223193756, 1, 0.3333333333333333, 0.2537453333333333, 0.761236, 0.1395119621451497, 0.1395119621451497, 18.38979837475822, 165.508185372824, 0.00693856915364691, 0.006938569153646908),
    @("M2", "Gnas", "Gcgr", "FAPs", 3, 1, 72.47344466666667, 217.420334, 0.04973458223193756, 0.04973458223193756, 3, 1, 1.194006666666667, 3.58202, 0.6564779367281227, 0.6564779367281226, 86.53377608829777, 778.80398479468, 0.03264965592765752, 0.03264965592765751),
    @("M2", "Gnas", "Gcgr", "M2", 3, 1, 72.47344466666667, 217.420334, 0.04973458223193756, 0.04973458223193756, 1, 0.3333333333333333, 0.1154886666666667, 0.346466, 0.06349693324616998, 0.06349693324616996, 8.369861493293778, 75.328753439644, 0.00315799344800749, 0.00315799344800749),
    @("M2", "Gnas", "Gcgr", "sCs", 3, 1, 72.47344466666667, 217.420334, 0.04973458223193756, 0.04973458223193756, 3, 1, 0.2555663333333333, 0.766699, 0.1405131678805576, 0.1405131678805576, 18.52177251749622, 166.695952657466, 0.00698836370262564, 0.006988363702625639),
    @("sCs", "Gnas", "Gcgr", "ECs", 3, 1, 343.5753833333333, 1030.72615, 0.2357770937183061, 0.235777093718306, 1, 0.3333333333333333, 0.2537453333333333, 0.761236, 0.1395119621451497, 0.1395119621451497, 87.18065016904444, 784.6258515214, 0.03289372497352174, 0.03289372497352173),
    @("sCs", "Gnas", "Gcgr", "FAPs", 3, 1, 343.5753833333333, 1030.72615, 0.2357770937183061, 0.235777093718306, 3, 1, 1.194006666666667, 3.58202, 0.6564779367281227, 0.6564779367281226, 410.2312982025555, 3692.081683823, 0.1547824600119468, 0.1547824600119468),
    @("sCs", "Gnas", "Gcgr", "M2", 3, 1, 343.5753833333333, 1030.72615, 0.2357770937183061, 0.235777093718306, 1, 0.3333333333333333, 0.1154886666666667, 0.346466, 0.06349693324616998, 0.06349693324616996, 39.67906292065555, 357.1115662859, 0.01497112238080724, 0.01497112238080724),
    @("sCs", "Gnas", "Gcgr", "sCs", 3, 1, 343.5753833333333, 1030.72615, 0.2357770937183061, 0.235777093718306, 3, 1, 0.2555663333333333, 0.766699, 0.1405131678805576, 0.1405131678805576, 87.80630094209444, 790.25670847885, 0.03312978635203031, 0.0331297863520303)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value = $rowData[$c]
    }
}
